$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 380 (the "「ここに見えるのは何？」" post) which causes all
# subsequent rows to shift up by one.
$ws.Rows.Item(380).Delete()
